$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.711.02"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").Value = "2.469.08"
$ws.Range("E3").Value = "  -0.87%  "

$ws.Range("E4").Value = "  +0.44%  "

$ws.Range("D5").Value = "'315.85"
$ws.Range("E5").Value = "  +0.73%  "

$ws.Range("D6").Value = "'92.82"
$ws.Range("E6").Value = "  -1.00%  "

$ws.Range("E7").Value = "  +0.44%  "

$ws.Range("E8").Value = "  +0.42%  "

$ws.Range("D9").Value = "'0.514"
$ws.Range("E9").Value = "  +3.04%  "

$ws.Range("D10").Value = "'32.65"
$ws.Range("E10").Value = "  -0.81%  "

$ws.Range("E11").Value = "  +6.38%  "

$ws.Range("E12").Value = "  +0.52%  "

$ws.Range("D13").Value = "2.847.02"
$ws.Range("E13").Value = "  -0.99%  "

$ws.Range("D14").Value = "'6.88"
$ws.Range("E14").Value = "  -0.13%  "

$ws.Range("D15").Value = "'15.86"
$ws.Range("E15").Value = "  +2.00%  "

$ws.Range("D16").Value = "2.442.75"
$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("D17").Value = "'0.781"
$ws.Range("E17").Value = "  +2.74%  "

$ws.Range("D18").Value = "41.668.95"
$ws.Range("E18").Value = "  -0.21%  "

$ws.Range("D19").Value = "'6.49"
$ws.Range("E19").Value = "  +2.35%  "

$ws.Range("E20").Value = "  +2.20%  "

$ws.Range("D21").Value = "'70.91"
$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("D22").Value = "'11.38"
$ws.Range("E22").Value = "  +1.05%  "

$ws.Range("D23").Value = "'239.17"
$ws.Range("E23").Value = "  +1.09%  "

$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("D25").Value = "'1.92"
$ws.Range("E25").Value = "  +0.64%  "

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").Value = "'24.56"
$ws.Range("E27").Value = "  -1.00%  "

$ws.Range("E28").Value = "  +0.68%  "

$ws.Range("D29").Value = "'9.77"
$ws.Range("E29").Value = "  +0.84%  "

$ws.Range("D30").Value = "'35.44"
$ws.Range("E30").Value = "  -2.11%  "

$ws.Range("D31").Value = "'155.56"
$ws.Range("E31").Value = "  +0.66%  "

$ws.Range("D32").Value = "'5.53"
$ws.Range("E32").Value = "  +2.09%  "

$ws.Range("E33").Value = "  +0.35%  "

$ws.Range("E34").Value = "  +0.72%  "

$ws.Range("D35").Value = "'2.48"
$ws.Range("E35").Value = "  -0.10%  "

$ws.Range("D36").Value = "'17.47"
$ws.Range("E36").Value = "  -5.01%  "

$ws.Range("D37").Value = "'2.90"
$ws.Range("E37").Value = "  -2.00%  "

$ws.Range("E38").Value = "  +1.17%  "

$ws.Range("E39").Value = "  +0.38%  "

$ws.Range("E40").Value = "  -1.83%  "

$ws.Range("D41").Value = "'3.96"
$ws.Range("E41").Value = "  -4.47%  "

$ws.Range("E42").Value = "  +0.35%  "

$ws.Range("D43").Value = "1.977.16"
$ws.Range("E43").Value = "  +1.42%  "

$ws.Range("D44").Value = "'19.03"
$ws.Range("E44").Value = "  -4.20%  "

$ws.Range("D45").Value = "'0.0283"
$ws.Range("E45").Value = "  -0.73%  "

$ws.Range("D46").Value = "'2.95"
$ws.Range("E46").Value = "  -1.17%  "

$ws.Range("D47").Value = "'9.03"
$ws.Range("E47").Value = "  +2.06%  "

$ws.Range("D48").Value = "2.702.99"
$ws.Range("E48").Value = "  -1.16%  "

$ws.Range("D49").Value = "'96.94"
$ws.Range("E49").Value = "  +0.39%  "

$ws.Range("D50").Value = "'67.15"
$ws.Range("E50").Value = "  -0.54%  "

$ws.Range("D51").Value = "'72.40"
$ws.Range("E51").Value = "  -1.29%  "

